$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "git branch"
$ws.Range("C9").Value = "لمعرفة التفرعات الموجودة وماهو التفرع الحالي"

$ws.Range("A10").Value = "git remote -v"
$ws.Range("C10").Value = "لمعرفة اسم remote"

$ws.Range("A11").Value = "git push origin master"
$ws.Range("C11").Value = "لرفع التغيرات من الفرع master في الجهاز المحلي(local) الى origin في الجهاز البعيد (remote)"

# Target stored column width is 69.42578125 characters; the engine quantizes
# ColumnWidth assignments to the nearest 1/6-character step, so 68.7 is the
# input that lands closest (stored width 69.5, the nearest reachable value).
$ws.Columns.Item(3).ColumnWidth = 68.7

$ws.Range("A12").Select()
